$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrigido montagem da EntradaDiretorio ---
# The "nome do arquivo" field of the directory entry was 20 bytes; it
# should be 16 bytes. This single value fix ripples through the A-column
# running "position" formulas below it.
$ws.Range("B18").Value = 16

# A new field, "Tamanho arquivo" (file size), was missing from the
# EntradaDiretorio layout; insert a row for it right after the
# "primeiro cluster" row (row 20), before "extensao".
$ws.Rows("21").Insert()

$ws.Range("A21").Formula = "=A20+B20"
$ws.Range("B21").Value = 4
$ws.Range("C21").Value = "Tamanho arquivo"

# Copy the formatting of the row above (same column styles used
# throughout this little table) onto the freshly inserted row.
$ws.Range("A20:C20").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)

# The rows that got pushed down keep referencing the literal cell
# addresses they had before the insert; re-point each running-position
# formula at the row now directly above it, matching the original
# "= previous row A + previous row B" pattern.
$ws.Range("A22").Formula = "=A21+B21"
$ws.Range("A23").Formula = "=A22+B22"
$ws.Range("A24").Formula = "=A23+B23"

# Page setup was switched to A4 portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection moved from the old E21 to C21 (inside the newly inserted /
# shifted rows).
$ws.Range("C21").Select()
